$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.408.04'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '3.015.27'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''508.61'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '''139.81'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '''7.57'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '''0.366'
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("D12").Value = '3.527.70'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '''26.35'
$ws.Range("E14").Value = '  +2.26%  '
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").Value = '57.392.21'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '''6.22'
$ws.Range("E17").Value = '  +4.12%  '
$ws.Range("D18").Value = '3.011.59'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = '''12.82'
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").Value = '''7.96'
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").Value = '''327.37'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("D24").Value = '''0.499'
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("D25").Value = '''64.56'
$ws.Range("E25").Value = '  +2.13%  '
$ws.Range("D26").Value = '''0.167'
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = '0.0₃0920'
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").Value = '''6.78'
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").Value = '''7.35'
$ws.Range("E30").Value = '  +4.26%  '
$ws.Range("D31").Value = '''1.81'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("D33").Value = '''20.58'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = '''4.78'
$ws.Range("E34").Value = '  +4.37%  '
$ws.Range("D35").Value = '''153.95'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '''5.89'
$ws.Range("D37").Value = '''1.28'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '''24.70'
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("D39").Value = '''0.0677'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").Value = '3.044.69'
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").Value = '''37.85'
$ws.Range("E41").Value = '  +2.52%  '
$ws.Range("E42").Value = '  +4.73%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''1.42'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '2.224.84'
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("D47").Value = '''0.981'
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").Value = '''6.06'
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("D49").Value = '''0.0239'
$ws.Range("E49").Value = '  -0.91%  '
$ws.Range("D50").Value = '''19.56'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  -5.22%  '
